$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Chad Oliver -> John Doe, with new diagnosis data)
$ws.Range("B2").Value = "John"
$ws.Range("C2").Value = "Doe"
$ws.Range("E2").Value = 111.2
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = "Chest Pain,Pressure,Fever,Tiredness"
$ws.Range("H2").Value = "Chronic Kidney Disease,Asthma,Dementia"
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = $false
$ws.Range("P2").Value = $false
$ws.Range("Q2").Value = "Very High Risk"

# Add new row 3 for Jane Doe
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Jane"
$ws.Range("C3").Value = "Doe"
$ws.Range("D3").Value = "chadoliver017@gmail.com"
$ws.Range("E3").Value = 113
$ws.Range("F3").Value = 68
$ws.Range("G3").Value = "Difficulty Breathing,Pressure,Fever,Pains,Headache,Loss of Taste,Discolouration"
$ws.Range("H3").Value = "Cancer,Cystic Fibrosis,Dementia,Down syndrome"
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = $false
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = "Very High Risk"
